$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I0 (I1) and IF (J1), copying the header style from H1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Populate I0 / IF data for rows 2-64
$data = @(
    @(2,6,6),
    @(3,11,11),
    @(4,4,5),
    @(5,8,9),
    @(6,7,7),
    @(7,5,6),
    @(8,7,7),
    @(9,4,5),
    @(10,7,7),
    @(11,5,5),
    @(12,6,6),
    @(13,7,7),
    @(14,8,8),
    @(15,8,8),
    @(16,4,5),
    @(17,7,7),
    @(18,7,7),
    @(19,7,7),
    @(20,7,7),
    @(21,6,6),
    @(22,7,7),
    @(23,7,7),
    @(24,8,8),
    @(25,6,6),
    @(26,7,7),
    @(27,7,7),
    @(28,7,7),
    @(29,6,6),
    @(30,6,6),
    @(31,6,7),
    @(32,7,7),
    @(33,7,7),
    @(34,7,7),
    @(35,7,7),
    @(36,10,10),
    @(37,7,7),
    @(38,7,8),
    @(39,9,9),
    @(40,6,6),
    @(41,8,8),
    @(42,7,7),
    @(43,6,6),
    @(44,6,7),
    @(45,7,7),
    @(46,6,6),
    @(47,7,8),
    @(48,8,9),
    @(49,7,7),
    @(50,10,11),
    @(51,6,6),
    @(52,9,9),
    @(53,7,8),
    @(54,5,6),
    @(55,7,7),
    @(56,7,8),
    @(57,9,9),
    @(58,2,3),
    @(59,8,9),
    @(60,7,7),
    @(61,6,6),
    @(62,4,4),
    @(63,4,4),
    @(64,4,4)
)

foreach ($row in $data) {
    $r = $row[0]
    $i = $row[1]
    $j = $row[2]
    $ws.Cells.Item($r, 9).Value = $i
    $ws.Cells.Item($r, 10).Value = $j
}

